$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear C2 (remove "Puvisuk")
$ws.Range("C2").Value = $null

# Row 3 becomes the former row 10 content (Stimmelmayr et al 2018)
$ws.Range("A3").Value = "Stimmelmayr et al 2018"
$ws.Range("B3").Value = "Stimmelmayr"
$ws.Range("C3").Value = "Utqiaġvik"

# Row 4 becomes the former row 11 content (Selendang)
$ws.Range("A4").Value = "Selendang"
$ws.Range("B4").Value = "Selendang"
$ws.Range("C4").Value = $null

# Delete old rows 5 through 11 which are no longer needed
$ws.Range("A5:F11").EntireRow.Delete()
